$wb = $excel.ActiveWorkbook

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2733.6667
$ws.Range("I43").Value = 1973.5
$ws.Range("J43").Value = 3113.75
$ws.Range("K43").Value = 1973.5
$ws.Range("L43").Value = 3113.75
$ws.Range("M43").Value = -1904.5
$ws.Range("N43").Value = -3251.75

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1297.1333
$ws.Range("I100").Value = 1408.4546
$ws.Range("K100").Value = 1408.4546
$ws.Range("M100").Value = -867.4546

# ALC row 108
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 49998.332
$ws.Range("J108").Value = 49998.332
$ws.Range("L108").Value = 49998.332
$ws.Range("N108").Value = -57678.332

# ALC row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 4945.5884
$ws.Range("J121").Value = 4945.5884
$ws.Range("L121").Value = 14836.7652
$ws.Range("N121").Value = -18330.7652

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6362818
$ws.Range("I137").Value = 298424.1
$ws.Range("J137").Value = 16672287
$ws.Range("K137").Value = 895272.2999999999
$ws.Range("L137").Value = 50016861
$ws.Range("M137").Value = -892722.2999999999
$ws.Range("N137").Value = -50021961

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6609.8057
$ws.Range("I138").Value = 2381.8462
$ws.Range("K138").Value = 7145.5386
$ws.Range("M138").Value = -2005.5386

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2326.4546
$ws.Range("I141").Value = 1899.1111
$ws.Range("K141").Value = 5697.3333
$ws.Range("M141").Value = -517.3333000000002

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14543.162
$ws.Range("I32").Value = 13158.849
$ws.Range("K32").Value = 13158.849
$ws.Range("M32").Value = -12871.849

# ARM row 44
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 60000
$ws.Range("J44").Value = 60000
$ws.Range("L44").Value = 60000
$ws.Range("N44").Value = -60976

# ARM row 55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 49333.332
$ws.Range("J55").Value = 49333.332
$ws.Range("L55").Value = 49333.332
$ws.Range("N55").Value = -49963.332

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4048
$ws.Range("I61").Value = 3679.8918
$ws.Range("K61").Value = 3679.8918
$ws.Range("M61").Value = -3467.8918

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1193.2693
$ws.Range("I74").Value = 1038.9474
$ws.Range("J74").Value = 1612.1428
$ws.Range("K74").Value = 1038.9474
$ws.Range("L74").Value = 1612.1428
$ws.Range("M74").Value = -164.9474
$ws.Range("N74").Value = -3360.1428

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1193.2693
$ws.Range("I77").Value = 1038.9474
$ws.Range("J77").Value = 1612.1428
$ws.Range("K77").Value = 5194.737
$ws.Range("L77").Value = 8060.714
$ws.Range("M77").Value = -826.7370000000001
$ws.Range("N77").Value = -16796.714

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3445.1428
$ws.Range("I122").Value = 2108.2
$ws.Range("K122").Value = 6324.599999999999
$ws.Range("M122").Value = -3874.599999999999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4048
$ws.Range("I136").Value = 3679.8918
$ws.Range("K136").Value = 11039.6754
$ws.Range("M136").Value = -8489.6754

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 93579.664
$ws.Range("J132").Value = 93579.664
$ws.Range("L132").Value = 93579.664
$ws.Range("N132").Value = -103699.664

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6062.7393
$ws.Range("I134").Value = 2869.087
$ws.Range("J134").Value = 9256.392
$ws.Range("K134").Value = 8607.261
$ws.Range("L134").Value = 27769.176
$ws.Range("M134").Value = -6072.261
$ws.Range("N134").Value = -32839.176

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 50002636
$ws.Range("I31").Value = 55557596
$ws.Range("K31").Value = 55557596
$ws.Range("M31").Value = -55557301

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 50002636
$ws.Range("I34").Value = 55557596
$ws.Range("K34").Value = 55557596
$ws.Range("M34").Value = -55557394

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2023.9032
$ws.Range("I58").Value = 1554.2
$ws.Range("K58").Value = 1554.2
$ws.Range("M58").Value = -1351.2

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6542.6816
$ws.Range("I99").Value = 6895.0415
$ws.Range("J99").Value = 6119.85
$ws.Range("K99").Value = 6895.0415
$ws.Range("L99").Value = 6119.85
$ws.Range("M99").Value = -5397.0415
$ws.Range("N99").Value = -9115.85

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5797.45
$ws.Range("I122").Value = 4226.846
$ws.Range("J122").Value = 8714.286
$ws.Range("K122").Value = 12680.538
$ws.Range("L122").Value = 26142.858
$ws.Range("M122").Value = -10230.538
$ws.Range("N122").Value = -31042.858

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6542.6816
$ws.Range("I126").Value = 6895.0415
$ws.Range("J126").Value = 6119.85
$ws.Range("K126").Value = 20685.1245
$ws.Range("L126").Value = 18359.55
$ws.Range("M126").Value = -18215.1245
$ws.Range("N126").Value = -23299.55

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2023.9032
$ws.Range("I136").Value = 1554.2
$ws.Range("K136").Value = 4662.6
$ws.Range("M136").Value = -2112.6

# CUL row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 257.29166
$ws.Range("I14").Value = 257.29166
$ws.Range("K14").Value = 771.8749799999999
$ws.Range("M14").Value = -598.8749799999999

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1717.8182
$ws.Range("I122").Value = 2352
$ws.Range("J122").Value = 1576.8889
$ws.Range("K122").Value = 21168
$ws.Range("L122").Value = 14192.0001
$ws.Range("M122").Value = -18718
$ws.Range("N122").Value = -19092.0001

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4984.75
$ws.Range("I70").Value = 4029
$ws.Range("J70").Value = 6152.8887
$ws.Range("K70").Value = 4029
$ws.Range("L70").Value = 6152.8887
$ws.Range("M70").Value = -3759
$ws.Range("N70").Value = -6692.8887

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4984.75
$ws.Range("I73").Value = 4029
$ws.Range("J73").Value = 6152.8887
$ws.Range("K73").Value = 4029
$ws.Range("L73").Value = 6152.8887
$ws.Range("M73").Value = -3093
$ws.Range("N73").Value = -8024.8887

# GSM row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# GSM row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 73047.5
$ws.Range("J140").Value = 96386
$ws.Range("L140").Value = 96386
$ws.Range("N140").Value = -106746

# GSM row 141
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 100428.5
$ws.Range("J141").Value = 100428.5
$ws.Range("L141").Value = 100428.5
$ws.Range("N141").Value = -110788.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1101.2
$ws.Range("I22").Value = 852.1111
$ws.Range("K22").Value = 852.1111
$ws.Range("M22").Value = -557.1111

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1101.2
$ws.Range("I27").Value = 852.1111
$ws.Range("K27").Value = 852.1111
$ws.Range("M27").Value = -745.1111

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 480.86667
$ws.Range("I55").Value = 111.5
$ws.Range("K55").Value = 111.5
$ws.Range("M55").Value = 61.5

# LTW row 137
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 90000
$ws.Range("J137").Value = 90000
$ws.Range("L137").Value = 90000
$ws.Range("N137").Value = -100200

# LTW row 139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

# LTW row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 130000
$ws.Range("J141").Value = 130000
$ws.Range("L141").Value = 130000
$ws.Range("N141").Value = -140360

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14665.625
$ws.Range("I54").Value = 6725
$ws.Range("K54").Value = 6725
$ws.Range("M54").Value = -6205

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4925
$ws.Range("I62").Value = 4850
$ws.Range("K62").Value = 4850
$ws.Range("M62").Value = -4226

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4925
$ws.Range("I65").Value = 4850
$ws.Range("K65").Value = 24250
$ws.Range("M65").Value = -21130

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6603.75
$ws.Range("I81").Value = 4473.4165
$ws.Range("J81").Value = 9799.25
$ws.Range("K81").Value = 8946.833000000001
$ws.Range("L81").Value = 19598.5
$ws.Range("M81").Value = -7885.833000000001
$ws.Range("N81").Value = -21720.5

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 6603.75
$ws.Range("I84").Value = 4473.4165
$ws.Range("J84").Value = 9799.25
$ws.Range("K84").Value = 44734.165
$ws.Range("L84").Value = 97992.5
$ws.Range("M84").Value = -39430.165
$ws.Range("N84").Value = -108600.5

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 856.2857
$ws.Range("I107").Value = 727.2857
$ws.Range("J107").Value = 985.2857
$ws.Range("K107").Value = 2181.8571
$ws.Range("L107").Value = 2955.8571
$ws.Range("M107").Value = -261.8571000000002
$ws.Range("N107").Value = -6795.8571

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3895.6667
$ws.Range("I122").Value = 2871.5557
$ws.Range("K122").Value = 8614.667099999999
$ws.Range("M122").Value = -6164.667099999999

# WVR row 124
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 39999
$ws.Range("J124").Value = 39999
$ws.Range("L124").Value = 39999
$ws.Range("N124").Value = -49819

# WVR row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 92500
$ws.Range("J137").Value = 92500
$ws.Range("L137").Value = 92500
